$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume (E) columns to text so numeric-looking values
# (e.g. "1.00", "0.608") are preserved exactly as strings, matching the
# original inline-string cell contents instead of being parsed as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "37.035.95"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").Value = "2.026.68"
$ws.Range("E3").Value = "  -2.02%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "226.63"
$ws.Range("E5").Value = "  -2.34%  "

$ws.Range("D6").Value = "0.608"
$ws.Range("E6").Value = "  -3.61%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "55.15"
$ws.Range("E8").Value = "  -3.88%  "

$ws.Range("D9").Value = "0.381"
$ws.Range("E9").Value = "  -2.16%  "

$ws.Range("D10").Value = "0.0792"
$ws.Range("E10").Value = "  +1.81%  "

$ws.Range("E11").Value = "  -3.69%  "

$ws.Range("D12").Value = "2.326.49"
$ws.Range("E12").Value = "  -2.02%  "

$ws.Range("D13").Value = "14.35"
$ws.Range("E13").Value = "  -3.36%  "

$ws.Range("D14").Value = "20.40"
$ws.Range("E14").Value = "  -2.26%  "

$ws.Range("E15").Value = "  -1.98%  "

$ws.Range("D16").Value = "5.15"
$ws.Range("E16").Value = "  -3.14%  "

$ws.Range("D17").Value = "2.024.99"
$ws.Range("E17").Value = "  -1.99%  "

$ws.Range("D18").Value = "36.986.82"
$ws.Range("E18").Value = "  -1.01%  "

$ws.Range("D19").Value = "6.25"
$ws.Range("E19").Value = "  +5.25%  "

$ws.Range("D20").Value = "68.81"
$ws.Range("E20").Value = "  -2.36%  "

$ws.Range("D21").Value = "0.0₃0827"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").Value = "226.19"
$ws.Range("E22").Value = "  -0.65%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").Value = "2.42"
$ws.Range("E24").Value = "  +2.57%  "

$ws.Range("E25").Value = "  -7.27%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "9.27"
$ws.Range("E26").Value = "  -3.98%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "165.74"
$ws.Range("E27").Value = "  -2.23%  "

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.125"
$ws.Range("E28").Value = "  -6.11%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "18.74"
$ws.Range("E29").Value = "  -3.61%  "

$ws.Range("D30").Value = "1.34"
$ws.Range("E30").Value = "  -2.73%  "

$ws.Range("E31").Value = "  -4.63%  "

$ws.Range("D32").Value = "4.48"
$ws.Range("E32").Value = "  -2.65%  "

$ws.Range("D33").Value = "0.0619"
$ws.Range("E33").Value = "  -2.14%  "

$ws.Range("D34").Value = "4.45"
$ws.Range("E34").Value = "  -3.42%  "

$ws.Range("D35").Value = "2.36"
$ws.Range("E35").Value = "  -4.67%  "

$ws.Range("D36").Value = "1.84"
$ws.Range("E36").Value = "  +1.18%  "

$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.15%  "

$ws.Range("D38").Value = "3.17"
$ws.Range("E38").Value = "  -4.90%  "

$ws.Range("D39").Value = "5.30"
$ws.Range("E39").Value = "  +0.80%  "

$ws.Range("D40").Value = "17.39"
$ws.Range("E40").Value = "  +4.48%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0219"
$ws.Range("E41").Value = "  -4.79%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.491.79"
$ws.Range("E42").Value = "  +0.85%  "

$ws.Range("D43").Value = "95.61"
$ws.Range("E43").Value = "  -3.71%  "

$ws.Range("D44").Value = "0.0928"
$ws.Range("E44").Value = "  -2.87%  "

$ws.Range("D45").Value = "2.78"
$ws.Range("E45").Value = "  -4.39%  "

$ws.Range("D46").Value = "1.14"
$ws.Range("E46").Value = "  -5.14%  "

$ws.Range("D47").Value = "7.36"
$ws.Range("E47").Value = "  +1.35%  "

$ws.Range("D48").Value = "1.01"
$ws.Range("E48").Value = "  -3.10%  "

$ws.Range("E49").Value = "  -0.36%  "

$ws.Range("D50").Value = "2.213.30"
$ws.Range("E50").Value = "  -2.01%  "

$ws.Range("D51").Value = "3.64"
$ws.Range("E51").Value = "  -6.95%  "

# Restore the default (unstyled) cell style now that the text values are set,
# so the cells match the original workbook formatting (no explicit style index).
$ws.Range("D2:E51").Style = "Normal"